# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1765
    3  = 809
    6  = 39
    7  = 12175
    12 = 1124
    13 = 882
    14 = 13554
    15 = 13653
    20 = 1011
    23 = 2099
    24 = 194
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
